$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record")

$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

"done"
